# Finflux "error" scenario sheet: correct the lock-in-period warning copy
# (drop the stray "after") and make the error sheet the active tab/selection,
# matching the author's edit for this workbook.

$wb = $excel.ActiveWorkbook

$errSheet = $wb.Worksheets.Item("error")

# Make the "error" sheet the active sheet/tab (it becomes tabSelected, and
# the workbook's active tab moves from "RecurringDeposit Transaction" to it).
$errSheet.Activate()

# Fix the wording of the lock-in warning message shown in B1.
$errSheet.Range("B1").Value = "Withdrawals blocked until ``01 September 2015``."

# Move the sheet's selection onto the edited cell.
$errSheet.Range("B1").Select()
